$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H40").Value = 1644.1316
$ws.Range("I40").Value = 1495.409
$ws.Range("J40").Value = 1848.625
$ws.Range("K40").Value = 1495.409
$ws.Range("L40").Value = 1848.625
$ws.Range("M40").Value = -1320.409
$ws.Range("N40").Value = -2198.625
$ws.Range("H76").Value = 3484.2236
$ws.Range("I76").Value = 2995.5454
$ws.Range("J76").Value = 4156.1562
$ws.Range("K76").Value = 2995.5454
$ws.Range("L76").Value = 4156.1562
$ws.Range("M76").Value = -2680.5454
$ws.Range("N76").Value = -4786.1562
$ws.Range("H79").Value = 3484.2236
$ws.Range("I79").Value = 2995.5454
$ws.Range("J79").Value = 4156.1562
$ws.Range("K79").Value = 2995.5454
$ws.Range("L79").Value = 4156.1562
$ws.Range("M79").Value = -1903.5454
$ws.Range("N79").Value = -6340.1562
$ws.Range("H132").Value = 2530.7827
$ws.Range("I132").Value = 2373.0908
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 7119.2724
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -4589.2724
$ws.Range("N132").Value = -23060
$ws.Range("H138").Value = 3492.4707
$ws.Range("I138").Value = 756.5333000000001
$ws.Range("J138").Value = 5652.421
$ws.Range("K138").Value = 2269.5999
$ws.Range("L138").Value = 16957.263
$ws.Range("M138").Value = 2870.4001
$ws.Range("N138").Value = -27237.263

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 2731.72
$ws.Range("I45").Value = 2860
$ws.Range("K45").Value = 2860
$ws.Range("M45").Value = -2483
$ws.Range("H74").Value = 4267.606
$ws.Range("I74").Value = 5461
$ws.Range("J74").Value = 1880.8182
$ws.Range("K74").Value = 5461
$ws.Range("L74").Value = 1880.8182
$ws.Range("M74").Value = -4587
$ws.Range("N74").Value = -3628.8182
$ws.Range("H77").Value = 4267.606
$ws.Range("I77").Value = 5461
$ws.Range("J77").Value = 1880.8182
$ws.Range("K77").Value = 27305
$ws.Range("L77").Value = 9404.091
$ws.Range("M77").Value = -22937
$ws.Range("N77").Value = -18140.091
$ws.Range("H102").Value = 2954.2856
$ws.Range("I102").Value = 1995
$ws.Range("K102").Value = 1995
$ws.Range("M102").Value = -373
$ws.Range("H110").Value = 5716
$ws.Range("I110").Value = 5074
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 5074
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = -3029
$ws.Range("N110").Value = -11090
$ws.Range("H122").Value = 1015506.6
$ws.Range("I122").Value = 1046197.75
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 3138593.25
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -3136143.25
$ws.Range("N122").Value = -13000

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 2622.9033
$ws.Range("I86").Value = 2639.2
$ws.Range("J86").Value = 2593.2727
$ws.Range("K86").Value = 2639.2
$ws.Range("L86").Value = 2593.2727
$ws.Range("M86").Value = -1516.2
$ws.Range("N86").Value = -4839.2727
$ws.Range("H89").Value = 2622.9033
$ws.Range("I89").Value = 2639.2
$ws.Range("J89").Value = 2593.2727
$ws.Range("K89").Value = 13196
$ws.Range("L89").Value = 12966.3635
$ws.Range("M89").Value = -7580
$ws.Range("N89").Value = -24198.3635
$ws.Range("H122").Value = 44359.5
$ws.Range("J122").Value = 44359.5
$ws.Range("L122").Value = 44359.5
$ws.Range("N122").Value = -54159.5
$ws.Range("H134").Value = 1783.7333
$ws.Range("I134").Value = 1767.4814
$ws.Range("J134").Value = 1930
$ws.Range("K134").Value = 5302.4442
$ws.Range("L134").Value = 5790
$ws.Range("M134").Value = -2767.4442
$ws.Range("N134").Value = -10860

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2513.3257
$ws.Range("I31").Value = 1501.8438
$ws.Range("K31").Value = 1501.8438
$ws.Range("M31").Value = -1206.8438
$ws.Range("H34").Value = 2513.3257
$ws.Range("I34").Value = 1501.8438
$ws.Range("K34").Value = 1501.8438
$ws.Range("M34").Value = -1299.8438
$ws.Range("H58").Value = 2317
$ws.Range("I58").Value = 1798.2727
$ws.Range("J58").Value = 2755.923
$ws.Range("K58").Value = 1798.2727
$ws.Range("L58").Value = 2755.923
$ws.Range("M58").Value = -1595.2727
$ws.Range("N58").Value = -3161.923
$ws.Range("H132").Value = 1979.15
$ws.Range("I132").Value = 996.5833
$ws.Range("K132").Value = 2989.7499
$ws.Range("M132").Value = -459.7498999999998
$ws.Range("H136").Value = 2317
$ws.Range("I136").Value = 1798.2727
$ws.Range("J136").Value = 2755.923
$ws.Range("K136").Value = 5394.8181
$ws.Range("L136").Value = 8267.769
$ws.Range("M136").Value = -2844.8181
$ws.Range("N136").Value = -13367.769

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H82").Value = 20015
$ws.Range("J82").Value = 20015
$ws.Range("L82").Value = 60045
$ws.Range("N82").Value = -60857
$ws.Range("H85").Value = 20015
$ws.Range("J85").Value = 20015
$ws.Range("L85").Value = 60045
$ws.Range("N85").Value = -62853
$ws.Range("H97").Value = 980.8461
$ws.Range("I97").Value = 495.16666
$ws.Range("J97").Value = 1397.1428
$ws.Range("K97").Value = 1485.49998
$ws.Range("L97").Value = 4191.428400000001
$ws.Range("M97").Value = -989.4999800000001
$ws.Range("N97").Value = -5183.428400000001
$ws.Range("H131").Value = 742.7083
$ws.Range("J131").Value = 953.6129
$ws.Range("L131").Value = 2860.8387
$ws.Range("N131").Value = -12940.8387

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 3687.0967
$ws.Range("I80").Value = 2752.3809
$ws.Range("J80").Value = 5650
$ws.Range("K80").Value = 2752.3809
$ws.Range("L80").Value = 5650
$ws.Range("M80").Value = -1754.3809
$ws.Range("N80").Value = -7646
$ws.Range("H83").Value = 3687.0967
$ws.Range("I83").Value = 2752.3809
$ws.Range("J83").Value = 5650
$ws.Range("K83").Value = 13761.9045
$ws.Range("L83").Value = 28250
$ws.Range("M83").Value = -8769.904500000001
$ws.Range("N83").Value = -38234
$ws.Range("H102").Value = 1740.6666
$ws.Range("I102").Value = 1740.4138
$ws.Range("J102").Value = 1741.4
$ws.Range("K102").Value = 1740.4138
$ws.Range("L102").Value = 1741.4
$ws.Range("M102").Value = -118.4138
$ws.Range("N102").Value = -4985.4
$ws.Range("H107").Value = 308.89474
$ws.Range("I107").Value = 243.61539
$ws.Range("J107").Value = 450.33334
$ws.Range("K107").Value = 243.61539
$ws.Range("L107").Value = 450.33334
$ws.Range("M107").Value = 1676.38461
$ws.Range("N107").Value = -4290.33334
$ws.Range("H126").Value = 2622.5
$ws.Range("I126").Value = 1120.6666
$ws.Range("J126").Value = 3123.111
$ws.Range("K126").Value = 3361.9998
$ws.Range("L126").Value = 9369.332999999999
$ws.Range("M126").Value = -891.9998000000001
$ws.Range("N126").Value = -14309.333

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 2022.5
$ws.Range("I7").Value = 2083.4614
$ws.Range("K7").Value = 2083.4614
$ws.Range("M7").Value = -1971.4614
$ws.Range("H68").Value = 418212.4
$ws.Range("I68").Value = 2001499.6
$ws.Range("J68").Value = 1557.8948
$ws.Range("K68").Value = 2001499.6
$ws.Range("L68").Value = 1557.8948
$ws.Range("M68").Value = -2000750.6
$ws.Range("N68").Value = -3055.8948
$ws.Range("H71").Value = 418212.4
$ws.Range("I71").Value = 2001499.6
$ws.Range("J71").Value = 1557.8948
$ws.Range("K71").Value = 10007498
$ws.Range("L71").Value = 7789.474
$ws.Range("M71").Value = -10003754
$ws.Range("N71").Value = -15277.474
$ws.Range("H82").Value = 2037
$ws.Range("I82").Value = 1796.25
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1796.25
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1435.25
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2037
$ws.Range("I85").Value = 1796.25
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1796.25
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -548.25
$ws.Range("N85").Value = -5496
$ws.Range("H122").Value = 2060.4443
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2508.8
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 7526.400000000001
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -12426.4
$ws.Range("H126").Value = 2022.5
$ws.Range("I126").Value = 2083.4614
$ws.Range("K126").Value = 6250.3842
$ws.Range("M126").Value = -3780.3842
$ws.Range("H132").Value = 3212.9524
$ws.Range("I132").Value = 2866.4285
$ws.Range("J132").Value = 3906
$ws.Range("K132").Value = 8599.2855
$ws.Range("L132").Value = 11718
$ws.Range("M132").Value = -6069.2855
$ws.Range("N132").Value = -16778

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 782
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 980
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2940
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -7280
$ws.Range("H126").Value = 1217.6666
$ws.Range("I126").Value = 1159.5454
$ws.Range("J126").Value = 1377.5
$ws.Range("K126").Value = 3478.6362
$ws.Range("L126").Value = 4132.5
$ws.Range("M126").Value = -1008.6362
$ws.Range("N126").Value = -9072.5
$ws.Range("H132").Value = 2500.4092
$ws.Range("I132").Value = 1357.9286
$ws.Range("J132").Value = 4499.75
$ws.Range("K132").Value = 4073.7858
$ws.Range("L132").Value = 13499.25
$ws.Range("M132").Value = -1543.7858
$ws.Range("N132").Value = -18559.25
$ws.Range("H140").Value = 38817
$ws.Range("J140").Value = 38817
$ws.Range("L140").Value = 38817
$ws.Range("N140").Value = -49177
